$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 31.498119354248047
$ws.Range("C2").Value = 6.448275566101074
$ws.Range("D2").Value = 22.090225219726562
$ws.Range("E2").Value = 57.85714340209961
